$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SIQ")

# CYRS review points 9 to 16: rows 11-17 in the "Status" column (H) move
# from "Not answered" to "Answered".
$ws.Range("H11").Value = "Answered"
$ws.Range("H12").Value = "Answered"
$ws.Range("H13").Value = "Answered"
$ws.Range("H14").Value = "Answered"
$ws.Range("H15").Value = "Answered"
$ws.Range("H16").Value = "Answered"
$ws.Range("H17").Value = "Answered"

# Update the sheet view: scroll so row 14 is the top-left visible row/col E,
# and move the active selection to E16.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("E16").Select()
